$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D9").Value = "커뮤니티, 댓글 데이터로 데이터 사이언스 할 수 있다?"
$ws.Range("E9").Value = "https://pdsi.pabii.com/community-comments-not-for-ds/#utm_source=rss&utm_medium=rss&utm_campaign=community-comments-not-for-ds"

$ws.Range("D27").Value = "EMNLP 2022 Review"
$ws.Range("E27").Value = "https://tech.scatterlab.co.kr/emnlp2022-review/"

$ws.Range("D32").Value = "[Impala] with 문(clause) 결과셋을 임의 저장하지 않음"
$ws.Range("E32").Value = "https://dodonam.tistory.com/404"
